$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers: "level_kode" / "level_nama" -> "Kode Jenis Pengguna" / "Nama Jenis Pengguna"
$ws.Range("A1").Value = "Kode Jenis Pengguna"
$ws.Range("B1").Value = "Nama Jenis Pengguna"

# Widen the two data columns
$ws.Columns.Item(1).ColumnWidth = 20.75
$ws.Columns.Item(2).ColumnWidth = 19.75

# Move the active selection to C2 (was B8)
$ws.Range("C2").Select() | Out-Null

# Widen the workbook window (cosmetic bookViews setting)
$excel.ActiveWindow.Width = 18350
